$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 172.75
$ws.Range("I2").Value = 133.90909
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 133.90909
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -20.90908999999999
$ws.Range("N2").Value = -826
# Row 11
$ws.Range("H11").Value = 57.75
$ws.Range("I11").Value = 57.75
$ws.Range("K11").Value = 57.75
$ws.Range("M11").Value = 82.25
# Row 40
$ws.Range("H40").Value = 4649.9
$ws.Range("I40").Value = 4883.3335
$ws.Range("K40").Value = 4883.3335
$ws.Range("M40").Value = -4708.3335
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()
# Row 86
$ws.Range("H86").Value = 4666.3335
$ws.Range("I86").Value = 3749.5
$ws.Range("K86").Value = 3749.5
$ws.Range("M86").Value = -2626.5
# Row 89
$ws.Range("H89").Value = 4666.3335
$ws.Range("I89").Value = 3749.5
$ws.Range("K89").Value = 18747.5
$ws.Range("M89").Value = -13131.5
# Row 112
$ws.Range("H112").Value = 1524.88
$ws.Range("J112").Value = 1568.2084
$ws.Range("L112").Value = 4704.6252
$ws.Range("N112").Value = -6920.6252
# Row 137
$ws.Range("H137").Value = 2915.739
$ws.Range("I137").Value = 3098.7646
$ws.Range("K137").Value = 9296.293799999999
$ws.Range("M137").Value = -6746.293799999999
# Row 138
$ws.Range("H138").Value = 662.1429000000001
$ws.Range("I138").Value = 662.1429000000001
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 1986.4287
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 3153.5713
$ws.Range("N138").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2067530.1
$ws.Range("I32").Value = 952530.4399999999
$ws.Range("K32").Value = 952530.4399999999
$ws.Range("M32").Value = -952243.4399999999
# Row 45
$ws.Range("H45").Value = 107244.5
$ws.Range("I45").Value = 141993.33
$ws.Range("K45").Value = 141993.33
$ws.Range("M45").Value = -141616.33
# Row 80
$ws.Range("H80").Value = 60050
# Row 83
$ws.Range("H83").Value = 60050

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 827.1429000000001
$ws.Range("I22").Value = 827.1429000000001
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 827.1429000000001
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -654.1429000000001
$ws.Range("N22").ClearContents()
# Row 94
$ws.Range("H94").Value = 117653620
$ws.Range("I94").Value = 142864880
$ws.Range("K94").Value = 142864880
$ws.Range("M94").Value = -142864429
# Row 105
$ws.Range("H105").Value = 14445842
$ws.Range("I105").Value = 1429555.2
$ws.Range("J105").Value = 22728934
$ws.Range("K105").Value = 1429555.2
$ws.Range("L105").Value = 22728934
$ws.Range("M105").Value = -1427808.2
$ws.Range("N105").Value = -22732428
# Row 134
$ws.Range("H134").Value = 3692.5
$ws.Range("I134").Value = 4180
$ws.Range("K134").Value = 12540
$ws.Range("M134").Value = -10005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 10
$ws.Range("I17").Value = 10
$ws.Range("K17").Value = 10
$ws.Range("M17").Value = 164
# Row 31
$ws.Range("H31").Value = 2843699
$ws.Range("I31").Value = 2147.4688
$ws.Range("J31").Value = 10421170
$ws.Range("K31").Value = 2147.4688
$ws.Range("L31").Value = 10421170
$ws.Range("M31").Value = -1852.4688
$ws.Range("N31").Value = -10421760
# Row 34
$ws.Range("H34").Value = 2843699
$ws.Range("I34").Value = 2147.4688
$ws.Range("J34").Value = 10421170
$ws.Range("K34").Value = 2147.4688
$ws.Range("L34").Value = 10421170
$ws.Range("M34").Value = -1945.4688
$ws.Range("N34").Value = -10421574
# Row 52
$ws.Range("H52").Value = 27973.334
$ws.Range("J52").Value = 27973.334
$ws.Range("L52").Value = 27973.334
$ws.Range("N52").Value = -28561.334
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# Row 132
$ws.Range("H132").Value = 3751.1428
$ws.Range("I132").Value = 3111.64
$ws.Range("J132").Value = 5349.9
$ws.Range("K132").Value = 9334.92
$ws.Range("L132").Value = 16049.7
$ws.Range("M132").Value = -6804.92
$ws.Range("N132").Value = -21109.7
# Row 138
$ws.Range("H138").Value = 67844.766
$ws.Range("J138").Value = 67844.766
$ws.Range("L138").Value = 67844.766
$ws.Range("N138").Value = -78124.766

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 212.16667
$ws.Range("I13").Value = 212.16667
$ws.Range("K13").Value = 636.50001
$ws.Range("M13").Value = -468.50001
# Row 75
$ws.Range("H75").Value = 1533.3334
$ws.Range("I75").Value = 1166.6666
$ws.Range("J75").Value = 1900
$ws.Range("K75").Value = 3499.9998
$ws.Range("L75").Value = 5700
$ws.Range("M75").Value = -2501.9998
$ws.Range("N75").Value = -7696
# Row 78
$ws.Range("H78").Value = 1533.3334
$ws.Range("I78").Value = 1166.6666
$ws.Range("J78").Value = 1900
$ws.Range("K78").Value = 10499.9994
$ws.Range("L78").Value = 17100
$ws.Range("M78").Value = -5507.999400000001
$ws.Range("N78").Value = -27084
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Range("H41").Value = 1332.3334
# Row 80
$ws.Range("H80").Value = 90911100
$ws.Range("I80").Value = 200000860
$ws.Range("K80").Value = 200000860
$ws.Range("M80").Value = -199999862
# Row 83
$ws.Range("H83").Value = 90911100
$ws.Range("I83").Value = 200000860
$ws.Range("K83").Value = 1000004300
$ws.Range("M83").Value = -999999308
# Row 118
$ws.Range("H118").Value = 49997.5
$ws.Range("J118").Value = 49997.5
$ws.Range("L118").Value = 49997.5
$ws.Range("N118").Value = -53311.5
# Row 122
$ws.Range("H122").Value = 3851.182
$ws.Range("I122").Value = 3050.7778
$ws.Range("K122").Value = 9152.3334
$ws.Range("M122").Value = -6702.3334

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2409.6667
$ws.Range("I7").Value = 1902.7142
$ws.Range("J7").Value = 3119.4
$ws.Range("K7").Value = 1902.7142
$ws.Range("L7").Value = 3119.4
$ws.Range("M7").Value = -1790.7142
$ws.Range("N7").Value = -3343.4
# Row 126
$ws.Range("H126").Value = 2409.6667
$ws.Range("I126").Value = 1902.7142
$ws.Range("J126").Value = 3119.4
$ws.Range("K126").Value = 5708.142599999999
$ws.Range("L126").Value = 9358.200000000001
$ws.Range("M126").Value = -3238.142599999999
$ws.Range("N126").Value = -14298.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
# Row 132
$ws.Range("H132").Value = 3195.4666
$ws.Range("I132").Value = 3188.8215
$ws.Range("J132").Value = 3288.5
$ws.Range("K132").Value = 9566.4645
$ws.Range("L132").Value = 9865.5
$ws.Range("M132").Value = -7036.4645
$ws.Range("N132").Value = -14925.5
